# vendor login and forget password validations
#
# - TestData sheet: add a new "palepu" value in D1 (new column) and make
#   TestData the active/selected sheet (was SignUp).
# - SignUp sheet: no longer the selected/active sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # TestData
$ws2 = $wb.Worksheets.Item(2)   # SignUp

# New data point alongside email / browser / url rows.
$ws1.Range("D1").Value = "palepu"

# Make TestData the active sheet/tab (was SignUp before the edit) and
# select the full header row (A1:XFD1) as the current selection there.
$ws1.Activate()
$ws1.Rows.Item(1).Select()
